# Fruta / hortaliza, semanal
# This workbook ("Fruta, Vega Modelo de Temuco - Kiwi") receives its weekly
# refresh: two rows that used to sit at 607/608 (dated 45005 = 2023-03-20) are
# duplicated down into newly inserted rows 609/610, and the original rows
# 607/608 are updated in place with the new week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Duplicate the current rows 607:608 down by inserting two new blank rows
#    right after them, then pasting a copy of the (still untouched) 607:608
#    range into the freshly inserted slot - those are exactly the rows that
#    should reappear, unedited, at 609:610.
$ws.Rows("609:610").Insert()
$ws.Range("A607:T608").Copy()
$ws.Range("A609").PasteSpecial()
$excel.CutCopyMode = 0

# 2) Now update row 607 in place with this week's values.
$ws.Range("D607").Value2 = 45015
$ws.Range("L607").Value2 = "Especial"
$ws.Range("M607").Value2 = 155

# 3) Update row 608 in place with this week's values.
$ws.Range("D608").Value2 = 45015
$ws.Range("M608").Value2 = 235
$ws.Range("N608").Value2 = 18000
$ws.Range("P608").Value2 = 18936
$ws.Range("S608").Value2 = 1052
